$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 435.4
$ws.Range("J2").Value = 277
$ws.Range("L2").Value = 277
$ws.Range("N2").Value = -503
$ws.Range("H5").Value = 343.33334
$ws.Range("I5").Value = 390
$ws.Range("K5").Value = 390
$ws.Range("M5").Value = -275
$ws.Range("H51").Value = 5188
$ws.Range("J51").Value = 5365
$ws.Range("L51").Value = 5365
$ws.Range("N51").Value = -6333
$ws.Range("H58").Value = 1167.4546
$ws.Range("I58").Value = 307.75
$ws.Range("K58").Value = 923.25
$ws.Range("M58").Value = -773.25
$ws.Range("H74").Value = 4544.625
$ws.Range("I74").Value = 4890.75
$ws.Range("J74").Value = 4198.5
$ws.Range("K74").Value = 4890.75
$ws.Range("L74").Value = 4198.5
$ws.Range("M74").Value = -3954.75
$ws.Range("N74").Value = -6070.5
$ws.Range("H76").Value = 7812073
$ws.Range("I76").Value = 11713312
$ws.Range("J76").Value = 9595
$ws.Range("K76").Value = 11713312
$ws.Range("L76").Value = 9595
$ws.Range("M76").Value = -11712997
$ws.Range("N76").Value = -10225
$ws.Range("H77").Value = 4544.625
$ws.Range("I77").Value = 4890.75
$ws.Range("J77").Value = 4198.5
$ws.Range("K77").Value = 24453.75
$ws.Range("L77").Value = 20992.5
$ws.Range("M77").Value = -19773.75
$ws.Range("N77").Value = -30352.5
$ws.Range("H79").Value = 7812073
$ws.Range("I79").Value = 11713312
$ws.Range("J79").Value = 9595
$ws.Range("K79").Value = 11713312
$ws.Range("L79").Value = 9595
$ws.Range("M79").Value = -11712220
$ws.Range("N79").Value = -11779
$ws.Range("H94").Value = 2176.2222
$ws.Range("I94").Value = 1898.7142
$ws.Range("K94").Value = 1898.7142
$ws.Range("M94").Value = -1447.7142
$ws.Range("H131").Value = 1615.8823
$ws.Range("I131").Value = 651.63635
$ws.Range("J131").Value = 3383.6667
$ws.Range("K131").Value = 1954.90905
$ws.Range("L131").Value = 10151.0001
$ws.Range("M131").Value = 3085.09095
$ws.Range("N131").Value = -20231.0001
$ws.Range("H132").Value = 995.2683
$ws.Range("I132").Value = 1007.1795
$ws.Range("K132").Value = 3021.5385
$ws.Range("M132").Value = -491.5384999999997
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4847.256
$ws.Range("I32").Value = 3407.8064
$ws.Range("J32").Value = 8565.833000000001
$ws.Range("K32").Value = 3407.8064
$ws.Range("L32").Value = 8565.833000000001
$ws.Range("M32").Value = -3120.8064
$ws.Range("N32").Value = -9139.833000000001
$ws.Range("H61").Value = 2265.5386
$ws.Range("I61").Value = 1372.1305
$ws.Range("K61").Value = 1372.1305
$ws.Range("M61").Value = -1160.1305
$ws.Range("H122").Value = 1772.5
$ws.Range("I122").Value = 1848.1875
$ws.Range("K122").Value = 5544.5625
$ws.Range("M122").Value = -3094.5625
$ws.Range("H136").Value = 2265.5386
$ws.Range("I136").Value = 1372.1305
$ws.Range("K136").Value = 4116.3915
$ws.Range("M136").Value = -1566.3915

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 49999
$ws.Range("J76").Value = 49999
$ws.Range("L76").Value = 49999
$ws.Range("N76").Value = -50629
$ws.Range("H79").Value = 49999
$ws.Range("J79").Value = 49999
$ws.Range("L79").Value = 49999
$ws.Range("N79").Value = -52183
$ws.Range("H80").Value = 7757.3076
$ws.Range("I80").Value = 33.25
$ws.Range("J80").Value = 11190.223
$ws.Range("K80").Value = 33.25
$ws.Range("L80").Value = 11190.223
$ws.Range("M80").Value = 964.75
$ws.Range("N80").Value = -13186.223
$ws.Range("H83").Value = 7757.3076
$ws.Range("I83").Value = 33.25
$ws.Range("J83").Value = 11190.223
$ws.Range("K83").Value = 166.25
$ws.Range("L83").Value = 55951.115
$ws.Range("M83").Value = 4825.75
$ws.Range("N83").Value = -65935.11499999999
$ws.Range("H86").Value = 168287.83
$ws.Range("I86").Value = 1819.8
$ws.Range("K86").Value = 1819.8
$ws.Range("M86").Value = -696.8
$ws.Range("H89").Value = 168287.83
$ws.Range("I89").Value = 1819.8
$ws.Range("K89").Value = 9099
$ws.Range("M89").Value = -3483
$ws.Range("H105").Value = 2500
$ws.Range("I105").Value = 2500
$ws.Range("K105").Value = 2500
$ws.Range("M105").Value = -753
$ws.Range("H107").Value = 2787.8
$ws.Range("I107").Value = 2487.5
$ws.Range("J107").Value = 3989
$ws.Range("K107").Value = 2487.5
$ws.Range("L107").Value = 3989
$ws.Range("M107").Value = -567.5
$ws.Range("N107").Value = -7829
$ws.Range("H135").Value = 31666.334
$ws.Range("J135").Value = 29999.5
$ws.Range("L135").Value = 29999.5
$ws.Range("N135").Value = -40139.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 29950
$ws.Range("J59").Value = 29950
$ws.Range("L59").Value = 29950
$ws.Range("N59").Value = -32240
$ws.Range("H107").Value = 355
$ws.Range("I107").Value = 293.17392
$ws.Range("K107").Value = 293.17392
$ws.Range("M107").Value = 1626.82608
$ws.Range("H132").Value = 1729.0344
$ws.Range("I132").Value = 1045.2
$ws.Range("K132").Value = 3135.6
$ws.Range("M132").Value = -605.6000000000004
$ws.Range("H134").Value = 888.4211
$ws.Range("I134").Value = 773.9375
$ws.Range("K134").Value = 2321.8125
$ws.Range("M134").Value = 213.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1.3333334
$ws.Range("H39").Value = 2900
$ws.Range("J39").Value = 2900
$ws.Range("L39").Value = 8700
$ws.Range("N39").Value = -9288
$ws.Range("H56").Value = 9687.875
$ws.Range("I56").Value = 9687.875
$ws.Range("K56").Value = 9687.875
$ws.Range("M56").Value = -9157.875
$ws.Range("H87").Value = 15390.4
$ws.Range("I87").Value = 976
$ws.Range("K87").Value = 2928
$ws.Range("M87").Value = -1680
$ws.Range("H90").Value = 15390.4
$ws.Range("I90").Value = 976
$ws.Range("K90").Value = 8784
$ws.Range("M90").Value = -2544
$ws.Range("H132").Value = 1744.3334
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 27000
$ws.Range("N132").Value = -32060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5571
$ws.Range("I70").Value = 5999.4
$ws.Range("K70").Value = 5999.4
$ws.Range("M70").Value = -5729.4
$ws.Range("H73").Value = 5571
$ws.Range("I73").Value = 5999.4
$ws.Range("K73").Value = 5999.4
$ws.Range("M73").Value = -5063.4
$ws.Range("H80").Value = 2444
$ws.Range("I80").Value = 2326.4
$ws.Range("K80").Value = 2326.4
$ws.Range("M80").Value = -1328.4
$ws.Range("H83").Value = 2444
$ws.Range("I83").Value = 2326.4
$ws.Range("K83").Value = 11632
$ws.Range("M83").Value = -6640
$ws.Range("H97").Value = 1095
$ws.Range("I97").Value = 556.25
$ws.Range("J97").Value = 3250
$ws.Range("K97").Value = 556.25
$ws.Range("L97").Value = 3250
$ws.Range("M97").Value = -60.25
$ws.Range("N97").Value = -4242
$ws.Range("H102").Value = 2258.8076
$ws.Range("I102").Value = 2302.3333
$ws.Range("K102").Value = 2302.3333
$ws.Range("M102").Value = -680.3332999999998
$ws.Range("H113").Value = 1621.4445
$ws.Range("I113").Value = 1333.25
$ws.Range("J113").Value = 1852
$ws.Range("K113").Value = 1333.25
$ws.Range("L113").Value = 1852
$ws.Range("M113").Value = 836.75
$ws.Range("N113").Value = -6192
$ws.Range("H126").Value = 65417.25
$ws.Range("J126").Value = 145452
$ws.Range("L126").Value = 436356
$ws.Range("N126").Value = -441296
$ws.Range("H132").Value = 4065.9707
$ws.Range("I132").Value = 3176
$ws.Range("K132").Value = 9528
$ws.Range("M132").Value = -6998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 258.30768
$ws.Range("I55").Value = 260
$ws.Range("K55").Value = 260
$ws.Range("M55").Value = -87
$ws.Range("H136").Value = 3997.8235
$ws.Range("I136").Value = 2211.875
$ws.Range("J136").Value = 5585.3335
$ws.Range("K136").Value = 6635.625
$ws.Range("L136").Value = 16756.0005
$ws.Range("M136").Value = -4085.625
$ws.Range("N136").Value = -21856.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 425
$ws.Range("I81").Value = 425
$ws.Range("K81").Value = 850
$ws.Range("M81").Value = 211
$ws.Range("H84").Value = 425
$ws.Range("I84").Value = 425
$ws.Range("K84").Value = 4250
$ws.Range("M84").Value = 1054
$ws.Range("H132").Value = 3338.111
$ws.Range("I132").Value = 3130.5
$ws.Range("K132").Value = 9391.5
$ws.Range("M132").Value = -6861.5
